$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Ontwerp van de applicatie" line: quantity dropped from 216 to 165 hours.
# E21 (=C21*D21) recalculates automatically.
$ws.Range("C21").Value = 165

# Subtotal now sums the whole line-item block (E20:E25) instead of just
# the first two rows, so future rows are included automatically.
# E28 (BTW) and E30 (Totaal) recalculate automatically from E26.
$ws.Range("E26").Formula = "=SUM(E20:E25)"

# Update the window scroll position / selection left behind when the
# workbook was last saved.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E27").Select()
